$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 3.5697000000000001
$ws.Range("F1").Value = 98.832300000000004
$ws.Range("E3").Value = 3.6168

$ws.Range("E4").Select()
